$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    7  = @(0.1838799416186027, 0.7695350798161764, 1.48474280445927, 1.218500227517119, 1.220715086378263, 38)
    8  = @(0.2202601043985484, 0.7025675926365696, 1.450233190988454, 1.204256281274237, 1.200273013662533, 37)
    9  = @(0.05458795430075868, 0.6814427688448653, 0.8599296228729056, 0.9273239039693227, 0.9497643907933797, 20)
    10 = @(0.3270187541837685, 0.5806918712726016, 0.7433548639178584, 0.8621802966420994, 0.8303300939129272, 13)
    11 = @(0.2613329855841933, 0.4854844728135347, 0.3301213188269617, 0.5745618494356911, 0.5720865204151965, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = 2 + $i
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
